$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update dialogue text (column B) that was rewritten in this commit ---
$ws.Range("B2").Value  = "Aren’t you very familiar with the layout of the manor?"
$ws.Range("B3").Value  = "As the butler, that’s part of my duty."
$ws.Range("B5").Value  = "When was the last time you saw the Lord today?"
$ws.Range("B6").Value  = "The last time I saw the Lord was shortly before the banquet was scheduled to begin."
$ws.Range("B8").Value  = " <color=#00CC00>(I remember you came to inform us of the original start time, which was——)</color>"
$ws.Range("B9").Value  = "Around 7 PM."
$ws.Range("B10").Value = "But the actual banquet started at 8 PM."
$ws.Range("B12").Value = "Roughly 6.30 PM."
$ws.Range("B15").Value = "I was there with Ling and Chen the entire time——they can confirm that."
$ws.Range("B17").Value = "So I went to the backyard right before 6.30 PM."
$ws.Range("B22").Value = "Although your statement is brief, there are a few key points to note——"
$ws.Range("B23").Value = " <color=#00CC00>(First, if Butler He said is true, then the Lord was still alive at 6.30 PM, and he left the backyard of his own accord.)</color>"
$ws.Range("B32").Value = " <color=#00CC00>(If the Lord was killed after it began raining, the killer’s clothes would probably have gotten wet, right?)</color>"
$ws.Range("B34").Value = "Are the clothes on the Lord’s body now the same ones you saw him wearing last time?"
$ws.Range("B40").Value = "After the 7 PM, I spent most of my time in the banquet hall with everyone else, waiting for the Lord to show up."
$ws.Range("B43").Value = "No. Before I helped the Young Lord search for the Lord, I only moved between the banquet hall and the kitchen."

# --- Row height adjustments ---
$ws.Rows.Item(2).RowHeight  = 17
$ws.Rows.Item(17).RowHeight = 17
$ws.Rows.Item(23).RowHeight = 51
$ws.Rows.Item(40).RowHeight = 34
$ws.Rows.Item(43).RowHeight = 34

# --- Selection / view state ---
$ws.Range("B49:B50").Select()
